$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 1.114 TL"
$ws.Range("D14").Value = ""
